# Update cryptos list prices and 1h volume percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.818.56"
$ws.Range("D3").Value = "1.624.94"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "214.42"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "0.500"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").Value = "0.0630"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").Value = "19.62"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").Value = "1.848.46"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "4.24"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "1.609.31"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("E15").Value = "  -2.89%  "
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").Value = "62.49"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "25.807.82"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Value = "1.00"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "192.36"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("D22").Value = "9.93"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "6.22"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "141.65"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("D29").Value = "15.43"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "1.24"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").Value = "0.0496"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").Value = "3.21"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "2.39"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("D36").Value = "0.900"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "1.125.39"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").Value = "0.545"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").Value = "2.47"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "99.25"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").Value = "1.759.65"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").Value = "0.0₆0111"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("E48").Value = "  +3.39%  "
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").Value = "7.56"
$ws.Range("E51").Value = "  +1.08%  "
